# Insert a new weekly price record for "Terminal La Palmera de La Serena - Cebollín"
# The new record is inserted as row 116, pushing the existing rows 116-223 down to 117-224.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 116 (shifts rows 116..223 down to 117..224,
# carrying formatting such as the date number format on column D along with them).
$ws.Rows(116).Insert()

# Populate the newly inserted row 116 with the new data point.
$ws.Range("A116").Value = 8
$ws.Range("B116").Value = "Terminal La Palmera de La Serena"
$ws.Range("C116").Value = "Coquimbo"
$ws.Range("D116").Value = 44778
$ws.Range("E116").Value = 4
$ws.Range("F116").Value = 100112037
$ws.Range("G116").Value = 'Cebollín'
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 2600
$ws.Range("K116").Value = 1400
$ws.Range("L116").Value = 1600
$ws.Range("M116").Value = 1500
$ws.Range("N116").Value = '$/paquete 6 unidades'
$ws.Range("O116").Value = 'Provincia del Elquí'
$ws.Range("P116").Value = 250
$ws.Range("Q116").Value = 6
$ws.Range("R116").Value = "Hortaliza"
